# Corrected Figure 3 experimental data sheet names
#
# Renames the two data sheets to match the corrected figure panel
# references and clears the now-superfluous bold/border style that had
# been applied to the "Note" header cell (D1) on both sheets, restoring
# it to the workbook's default (unstyled) cell format.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Fig 3C"
$ws1.Range("D1").Style = "Normal"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Fig 3E"
$ws2.Range("D1").Style = "Normal"
